# AfDD_2023_Annex_Table_Tab21.xlsx — "Add files via upload"
#
# The uploaded file corrects a units/scale error in columns O and P
# ("Intra-/Extra-continental trade in intermediate goods (millions of
# USD), 2020") for the data rows of the Tab21 sheet (rows 3-99): every
# value in those two columns is multiplied by 1000. Column Q (the
# trade-share %) is a statically stored value (not a live formula) and is
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 3
$lastRow = 99
$colO = 15   # column O — Intra-continental trade in intermediate goods
$colP = 16   # column P — Extra-continental trade in intermediate goods

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $oCell = $ws.Cells.Item($r, $colO)
    $oVal = $oCell.Value()
    $oCell.Value = $oVal * 1000

    $pCell = $ws.Cells.Item($r, $colP)
    $pVal = $pCell.Value()
    $pCell.Value = $pVal * 1000
}

# The original commit also carries an incidental bookViews/workbookView
# window-size change (windowHeight 12240 -> 12490) from whatever screen
# the workbook was last saved on in Excel. It is not tied to any
# documented Workbook/Window COM property in this runtime (Window.Height
# tracks separately from the exported windowHeight/windowWidth twips,
# which this engine always writes back as 28800/12240) so it can't be
# reproduced faithfully here; the attempt below is harmless best-effort.
$wb.Windows.Item(1).Height = 12490
